$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (item #1)
$ws.Range("A3").Value = 1
$ws.Range("A3").NumberFormat = "0"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "11-12-2025"
$ws.Range("C3").Value = "Shahul hameed"
$ws.Range("D3").Value = 9946508158
$ws.Range("D3").NumberFormat = "0"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "27-12-2025"
$ws.Range("F3").Value = "Mohamed Aslam A S"
$ws.Range("G3").Value = "Loss"
$ws.Range("H3").Value = "PRODUCT"
$ws.Range("I3").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "kids suit lavender"

# Row 4 (item #2)
$ws.Range("A4").Value = 2
$ws.Range("A4").NumberFormat = "0"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "11-12-2025"
$ws.Range("C4").Value = "salman"
$ws.Range("D4").Value = 7736144146
$ws.Range("D4").NumberFormat = "0"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "28-12-2025"
$ws.Range("F4").Value = "MAHESH C"
$ws.Range("G4").Value = "Loss"
$ws.Range("H4").Value = "PRODUCT"
$ws.Range("I4").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J4").Value = "-"
$ws.Range("K4").Value = "embro kurta"

# Row 5 (item #3)
$ws.Range("A5").Value = 3
$ws.Range("A5").NumberFormat = "0"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "11-12-2025"
$ws.Range("C5").Value = "Raziq"
$ws.Range("D5").Value = 7034581604
$ws.Range("D5").NumberFormat = "0"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "11-01-2026"
$ws.Range("F5").Value = "MUHAMMED ASLAM VB"
$ws.Range("G5").Value = "Loss"
$ws.Range("H5").Value = "PRODUCT"
$ws.Range("I5").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J5").Value = "-"
$ws.Range("K5").Value = "."

# Row 6 (item #4)
$ws.Range("A6").Value = 4
$ws.Range("A6").NumberFormat = "0"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "11-12-2025"
$ws.Range("C6").Value = "Yaseen"
$ws.Range("D6").Value = 9539373567
$ws.Range("D6").NumberFormat = "0"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "03-01-2026"
$ws.Range("F6").Value = "RASAL"
$ws.Range("G6").Value = "Loss"
$ws.Range("H6").Value = "PRODUCT"
$ws.Range("I6").Value = "REQUIRED DESIGN NOT AVAILABLE"
$ws.Range("J6").Value = "-"
$ws.Range("K6").Value = "DOUBLE BREAST"

# Row 7 (item #5)
$ws.Range("A7").Value = 5
$ws.Range("A7").NumberFormat = "0"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "12-12-2025"
$ws.Range("C7").Value = "Edwin"
$ws.Range("D7").Value = 9061310529
$ws.Range("D7").NumberFormat = "0"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "27-12-2025"
$ws.Range("F7").Value = "MAHESH C"
$ws.Range("G7").Value = "Loss"
$ws.Range("H7").Value = "PRODUCT"
$ws.Range("I7").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = "embro kurta"

# Row 8 (item #6)
$ws.Range("A8").Value = 6
$ws.Range("A8").NumberFormat = "0"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "13-12-2025"
$ws.Range("C8").Value = "Raeed"
$ws.Range("D8").Value = 9207897660
$ws.Range("D8").NumberFormat = "0"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "19-12-2025"
$ws.Range("F8").Value = "RASAL"
$ws.Range("G8").Value = "Loss"
$ws.Range("H8").Value = "PRODUCT"
$ws.Range("I8").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J8").Value = "-"
$ws.Range("K8").Value = "BOOTCUT PANTS"

# Row 9 (item #7)
$ws.Range("A9").Value = 7
$ws.Range("A9").NumberFormat = "0"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "14-12-2025"
$ws.Range("C9").Value = "zain"
$ws.Range("D9").Value = 7510601421
$ws.Range("D9").NumberFormat = "0"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "24-12-2025"
$ws.Range("F9").Value = "RASAL"
$ws.Range("G9").Value = "Loss"
$ws.Range("H9").Value = "PRODUCT"
$ws.Range("I9").Value = "REQUIRED DESIGN NOT AVAILABLE"
$ws.Range("J9").Value = "-"
$ws.Range("K9").Value = "lose"

# Row 10 (item #8)
$ws.Range("A10").Value = 8
$ws.Range("A10").NumberFormat = "0"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "15-12-2025"
$ws.Range("C10").Value = "NIMAL VB"
$ws.Range("D10").Value = 7025830557
$ws.Range("D10").NumberFormat = "0"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "23-12-2025"
$ws.Range("F10").Value = "MUHAMMED ASLAM VB"
$ws.Range("G10").Value = "Loss"
$ws.Range("H10").Value = "PRODUCT"
$ws.Range("I10").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J10").Value = "-"
$ws.Range("K10").Value = "CLOSED INDOWESTERN"

# Row 11 (item #9)
$ws.Range("A11").Value = 9
$ws.Range("A11").NumberFormat = "0"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "15-12-2025"
$ws.Range("C11").Value = "muzamil"
$ws.Range("D11").Value = 9037718484
$ws.Range("D11").NumberFormat = "0"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "19-12-2025"
$ws.Range("F11").Value = "MAHESH C"
$ws.Range("G11").Value = "Loss"
$ws.Range("H11").Value = "PRICING"
$ws.Range("I11").Value = "RENT TO HIGH"
$ws.Range("J11").Value = "-"
$ws.Range("K11").Value = "loss"

# Row 12 (item #10)
$ws.Range("A12").Value = 10
$ws.Range("A12").NumberFormat = "0"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "15-12-2025"
$ws.Range("C12").Value = "ajmal"
$ws.Range("D12").Value = 7994281115
$ws.Range("D12").NumberFormat = "0"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "21-12-2025"
$ws.Range("F12").Value = "RASAL"
$ws.Range("G12").Value = "Loss"
$ws.Range("H12").Value = "PRODUCT"
$ws.Range("I12").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J12").Value = "-"
$ws.Range("K12").Value = "double breasted suit"

# Row 13 (item #11)
$ws.Range("A13").Value = 11
$ws.Range("A13").NumberFormat = "0"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "16-12-2025"
$ws.Range("C13").Value = "Abdul razzaq"
$ws.Range("D13").Value = 9746697775
$ws.Range("D13").NumberFormat = "0"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "17-12-2025"
$ws.Range("F13").Value = "MAHESH C"
$ws.Range("G13").Value = "Loss"
$ws.Range("H13").Value = "SIZE NOT SUITABLE"
$ws.Range("I13").Value = "SIZE TOO LARGE"
$ws.Range("J13").Value = "-"
$ws.Range("K13").Value = "BIG SIZE"

# Row 14 (item #12)
$ws.Range("A14").Value = 12
$ws.Range("A14").NumberFormat = "0"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "16-12-2025"
$ws.Range("C14").Value = "fasil"
$ws.Range("D14").Value = 9961122822
$ws.Range("D14").NumberFormat = "0"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "25-12-2025"
$ws.Range("F14").Value = "RASAL"
$ws.Range("G14").Value = "Loss"
$ws.Range("H14").Value = "ENQUIRY"
$ws.Range("I14").Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = "GROOM IS NOT AVIALABLE"

# Row 15 (item #13)
$ws.Range("A15").Value = 13
$ws.Range("A15").NumberFormat = "0"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "16-12-2025"
$ws.Range("C15").Value = "shaheer"
$ws.Range("D15").Value = 7559857541
$ws.Range("D15").NumberFormat = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "28-12-2025"
$ws.Range("F15").Value = "MAHESH C"
$ws.Range("G15").Value = "Loss"
$ws.Range("H15").Value = "ENQUIRY"
$ws.Range("I15").Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Range("J15").Value = "-"

# Row 16 (item #14)
$ws.Range("A16").Value = 14
$ws.Range("A16").NumberFormat = "0"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "16-12-2025"
$ws.Range("C16").Value = "ajmal"
$ws.Range("D16").Value = 8139089882
$ws.Range("D16").NumberFormat = "0"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "14-01-2026"
$ws.Range("F16").Value = "MAHESH C"
$ws.Range("G16").Value = "Loss"
$ws.Range("H16").Value = "CUSTOMER INTERNAL ISSUES"
$ws.Range("I16").Value = "BUDGET RESTRICTIONS"
$ws.Range("J16").Value = "-"
$ws.Range("K16").Value = "CONFIRM TMRW"

# Row 17 (item #15)
$ws.Range("A17").Value = 15
$ws.Range("A17").NumberFormat = "0"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "17-12-2025"
$ws.Range("C17").Value = "aslam"
$ws.Range("D17").Value = 8592989565
$ws.Range("D17").NumberFormat = "0"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "31-12-2025"
$ws.Range("F17").Value = "MUHAMMED ASLAM VB"
$ws.Range("G17").Value = "Loss"
$ws.Range("H17").Value = "PRODUCT"
$ws.Range("I17").Value = "-"
$ws.Range("J17").Value = "-"
$ws.Range("K17").Value = "embro kurta"

# Row 18 (item #16)
$ws.Range("A18").Value = 16
$ws.Range("A18").NumberFormat = "0"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "17-12-2025"
$ws.Range("C18").Value = "thaha"
$ws.Range("D18").Value = 9847692714
$ws.Range("D18").NumberFormat = "0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "22-12-2025"
$ws.Range("F18").Value = "MAHESH C"
$ws.Range("G18").Value = "Loss"
$ws.Range("H18").Value = "PRODUCT"
$ws.Range("I18").Value = "PRODUCT NOT AVAILABLE"
$ws.Range("J18").Value = "-"
$ws.Range("K18").Value = "embro kurta"

# Row 19 (item #17)
$ws.Range("A19").Value = 17
$ws.Range("A19").NumberFormat = "0"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "18-12-2025"
$ws.Range("C19").Value = "faisal"
$ws.Range("D19").Value = 8714151963
$ws.Range("D19").NumberFormat = "0"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "03-01-2026"
$ws.Range("F19").Value = "MAHESH C"
$ws.Range("G19").Value = "Loss"
$ws.Range("H19").Value = "SIZE NOT SUITABLE"
$ws.Range("I19").Value = "SIZE TOO LARGE"
$ws.Range("J19").Value = "-"

# Row 20 (item #18)
$ws.Range("A20").Value = 18
$ws.Range("A20").NumberFormat = "0"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "19-12-2025"
$ws.Range("C20").Value = "sajadh"
$ws.Range("D20").Value = 9072022771
$ws.Range("D20").NumberFormat = "0"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "28-12-2025"
$ws.Range("F20").Value = "MAHESH C"
$ws.Range("G20").Value = "Loss"
$ws.Range("H20").Value = "PRODUCT"
$ws.Range("I20").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J20").Value = "-"
$ws.Range("K20").Value = "WHITE BENGALA NOT AVIALABLE IN OUR STORE"

# Row 21 (item #19)
$ws.Range("A21").Value = 19
$ws.Range("A21").NumberFormat = "0"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "21-12-2025"
$ws.Range("C21").Value = "nandhu"
$ws.Range("D21").Value = 8921731953
$ws.Range("D21").NumberFormat = "0"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "28-12-2025"
$ws.Range("F21").Value = "ANAS M J"
$ws.Range("G21").Value = "Loss"
$ws.Range("H21").Value = "PRODUCT"
$ws.Range("I21").Value = "REQUIRED MODEL NOT AVAILABLE"
$ws.Range("J21").Value = "-"
$ws.Range("K21").Value = "double brust suit"

# Update column K width to 54 (characters value empirically tuned to yield OOXML width=54)
$ws.Columns.Item(11).ColumnWidth = 53.14

Write-Host "Edit complete"